$d = $word.ActiveDocument

# Locate the paragraph containing the "Ver no Jupiter ..." text and the
# paragraph containing the "(c) 2020 ..." copyright text, plus the
# (empty) paragraph that immediately precedes the first one. Together
# these three paragraphs are removed, while the surrounding paragraphs
# (the "LOQ4204: ..." requirement line before, and the blank / page-break
# paragraphs after) are left untouched.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*Ver no Jupiter*") {
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($txt -like "*Contact: luizeleno@usp.br*") {
        $endPara = $d.Paragraphs.Item($i)
    }
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
